$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 9.233028333333333 ; $ws.Cells.Item(2, 8).Value = 27.699085 ; $ws.Cells.Item(2, 9).Value = 0.2291653310312791 ; $ws.Cells.Item(2, 10).Value = 0.2338715303104729 ; $ws.Cells.Item(2, 13).Value = 2.906846333333333 ; $ws.Cells.Item(2, 14).Value = 8.720538999999999 ; $ws.Cells.Item(2, 15).Value = 0.005520525738044089 ; $ws.Cells.Item(2, 16).Value = 0.005624540846623205 ; $ws.Cells.Item(2, 17).Value = 26.83899455631277 ; $ws.Cells.Item(2, 18).Value = 241.550951006815 ; $ws.Cells.Item(2, 19).Value = 0.00126511310822557 ; $ws.Cells.Item(2, 20).Value = 0.001315419975093532
$ws.Cells.Item(3, 7).Value = 9.233028333333333 ; $ws.Cells.Item(3, 8).Value = 27.699085 ; $ws.Cells.Item(3, 9).Value = 0.2291653310312791 ; $ws.Cells.Item(3, 10).Value = 0.2338715303104729 ; $ws.Cells.Item(3, 13).Value = 185.8027443333333 ; $ws.Cells.Item(3, 14).Value = 557.408233 ; $ws.Cells.Item(3, 15).Value = 0.3528665483720876 ; $ws.Cells.Item(3, 16).Value = 0.3595150912979765 ; $ws.Cells.Item(3, 17).Value = 1715.522002840756 ; $ws.Cells.Item(3, 18).Value = 15439.69802556681 ; $ws.Cells.Item(3, 19).Value = 0.08086477936755433 ; $ws.Cells.Item(3, 20).Value = 0.08408034457156717
$ws.Cells.Item(4, 7).Value = 9.233028333333333 ; $ws.Cells.Item(4, 8).Value = 27.699085 ; $ws.Cells.Item(4, 9).Value = 0.2291653310312791 ; $ws.Cells.Item(4, 10).Value = 0.2338715303104729 ; $ws.Cells.Item(4, 13).Value = 137.0717086666666 ; $ws.Cells.Item(4, 14).Value = 411.2151259999999 ; $ws.Cells.Item(4, 15).Value = 0.2603191943704447 ; $ws.Cells.Item(4, 16).Value = 0.2652240042658267 ; $ws.Cells.Item(4, 17).Value = 1265.586969817745 ; $ws.Cells.Item(4, 18).Value = 11390.28272835971 ; $ws.Cells.Item(4, 19).Value = 0.05965613435169884 ; $ws.Cells.Item(4, 20).Value = 0.0620283437527203
$ws.Cells.Item(5, 7).Value = 9.233028333333333 ; $ws.Cells.Item(5, 8).Value = 27.699085 ; $ws.Cells.Item(5, 9).Value = 0.2291653310312791 ; $ws.Cells.Item(5, 10).Value = 0.2338715303104729 ; $ws.Cells.Item(5, 13).Value = 29.2127365 ; $ws.Cells.Item(5, 14).Value = 58.425473 ; $ws.Cells.Item(5, 15).Value = 0.05547925319534149 ; $ws.Cells.Item(5, 16).Value = 0.03768304451958546 ; $ws.Cells.Item(5, 17).Value = 269.7220237987008 ; $ws.Cells.Item(5, 18).Value = 1618.332142792205 ; $ws.Cells.Item(5, 19).Value = 0.01271392142387858 ; $ws.Cells.Item(5, 20).Value = 0.008812991288553132
$ws.Cells.Item(6, 7).Value = 9.233028333333333 ; $ws.Cells.Item(6, 8).Value = 27.699085 ; $ws.Cells.Item(6, 9).Value = 0.2291653310312791 ; $ws.Cells.Item(6, 10).Value = 0.2338715303104729 ; $ws.Cells.Item(6, 13).Value = 171.5584106666666 ; $ws.Cells.Item(6, 14).Value = 514.6752319999999 ; $ws.Cells.Item(6, 15).Value = 0.3258144783240821 ; $ws.Cells.Item(6, 16).Value = 0.331953319069988 ; $ws.Cells.Item(6, 17).Value = 1584.003666506969 ; $ws.Cells.Item(6, 18).Value = 14256.03299856272 ; $ws.Cells.Item(6, 19).Value = 0.0746653827799218 ; $ws.Cells.Item(6, 20).Value = 0.07763443072253881
$ws.Cells.Item(7, 7).Value = 19.39389166666666 ; $ws.Cells.Item(7, 8).Value = 58.181675 ; $ws.Cells.Item(7, 9).Value = 0.4813596843119293 ; $ws.Cells.Item(7, 10).Value = 0.4912450129048156 ; $ws.Cells.Item(7, 13).Value = 2.906846333333333 ; $ws.Cells.Item(7, 14).Value = 8.720538999999999 ; $ws.Cells.Item(7, 15).Value = 0.005520525738044089 ; $ws.Cells.Item(7, 16).Value = 0.005624540846623205 ; $ws.Cells.Item(7, 17).Value = 56.37506288031387 ; $ws.Cells.Item(7, 18).Value = 507.3755659228249 ; $ws.Cells.Item(7, 19).Value = 0.002657358526500783 ; $ws.Cells.Item(7, 20).Value = 0.002763027640783079
$ws.Cells.Item(8, 7).Value = 19.39389166666666 ; $ws.Cells.Item(8, 8).Value = 58.181675 ; $ws.Cells.Item(8, 9).Value = 0.4813596843119293 ; $ws.Cells.Item(8, 10).Value = 0.4912450129048156 ; $ws.Cells.Item(8, 13).Value = 185.8027443333333 ; $ws.Cells.Item(8, 14).Value = 557.408233 ; $ws.Cells.Item(8, 15).Value = 0.3528665483720876 ; $ws.Cells.Item(8, 16).Value = 0.3595150912979765 ; $ws.Cells.Item(8, 17).Value = 3603.43829497003 ; $ws.Cells.Item(8, 18).Value = 32430.94465473027 ; $ws.Cells.Item(8, 19).Value = 0.1698557303286282 ; $ws.Cells.Item(8, 20).Value = 0.1766099956641504
$ws.Cells.Item(9, 7).Value = 19.39389166666666 ; $ws.Cells.Item(9, 8).Value = 58.181675 ; $ws.Cells.Item(9, 9).Value = 0.4813596843119293 ; $ws.Cells.Item(9, 10).Value = 0.4912450129048156 ; $ws.Cells.Item(9, 13).Value = 137.0717086666666 ; $ws.Cells.Item(9, 14).Value = 411.2151259999999 ; $ws.Cells.Item(9, 15).Value = 0.2603191943704447 ; $ws.Cells.Item(9, 16).Value = 0.2652240042658267 ; $ws.Cells.Item(9, 17).Value = 2658.353868446227 ; $ws.Cells.Item(9, 18).Value = 23925.18481601604 ; $ws.Cells.Item(9, 19).Value = 0.125307165222493 ; $ws.Cells.Item(9, 20).Value = 0.1302899693982329
$ws.Cells.Item(10, 7).Value = 19.39389166666666 ; $ws.Cells.Item(10, 8).Value = 58.181675 ; $ws.Cells.Item(10, 9).Value = 0.4813596843119293 ; $ws.Cells.Item(10, 10).Value = 0.4912450129048156 ; $ws.Cells.Item(10, 13).Value = 29.2127365 ; $ws.Cells.Item(10, 14).Value = 58.425473 ; $ws.Cells.Item(10, 15).Value = 0.05547925319534149 ; $ws.Cells.Item(10, 16).Value = 0.03768304451958546 ; $ws.Cells.Item(10, 17).Value = 566.5486469678791 ; $ws.Cells.Item(10, 18).Value = 3399.291881807275 ; $ws.Cells.Item(10, 19).Value = 0.02670547580397117 ; $ws.Cells.Item(10, 20).Value = 0.0185116076913165
$ws.Cells.Item(11, 7).Value = 19.39389166666666 ; $ws.Cells.Item(11, 8).Value = 58.181675 ; $ws.Cells.Item(11, 9).Value = 0.4813596843119293 ; $ws.Cells.Item(11, 10).Value = 0.4912450129048156 ; $ws.Cells.Item(11, 13).Value = 171.5584106666666 ; $ws.Cells.Item(11, 14).Value = 514.6752319999999 ; $ws.Cells.Item(11, 15).Value = 0.3258144783240821 ; $ws.Cells.Item(11, 16).Value = 0.331953319069988 ; $ws.Cells.Item(11, 17).Value = 3327.185230974844 ; $ws.Cells.Item(11, 18).Value = 29944.6670787736 ; $ws.Cells.Item(11, 19).Value = 0.1568339544303361 ; $ws.Cells.Item(11, 20).Value = 0.1630704125103327
$ws.Cells.Item(12, 7).Value = 3.072172 ; $ws.Cells.Item(12, 8).Value = 9.216516 ; $ws.Cells.Item(12, 9).Value = 0.07625183070469947 ; $ws.Cells.Item(12, 10).Value = 0.07781775827797052 ; $ws.Cells.Item(12, 13).Value = 2.906846333333333 ; $ws.Cells.Item(12, 14).Value = 8.720538999999999 ; $ws.Cells.Item(12, 15).Value = 0.005520525738044089 ; $ws.Cells.Item(12, 16).Value = 0.005624540846623205 ; $ws.Cells.Item(12, 17).Value = 8.930331913569331 ; $ws.Cells.Item(12, 18).Value = 80.372987222124 ; $ws.Cells.Item(12, 19).Value = 0.0004209501939782739 ; $ws.Cells.Item(12, 20).Value = 0.0004376891600270963
$ws.Cells.Item(13, 7).Value = 3.072172 ; $ws.Cells.Item(13, 8).Value = 9.216516 ; $ws.Cells.Item(13, 9).Value = 0.07625183070469947 ; $ws.Cells.Item(13, 10).Value = 0.07781775827797052 ; $ws.Cells.Item(13, 13).Value = 185.8027443333333 ; $ws.Cells.Item(13, 14).Value = 557.408233 ; $ws.Cells.Item(13, 15).Value = 0.3528665483720876 ; $ws.Cells.Item(13, 16).Value = 0.3595150912979765 ; $ws.Cells.Item(13, 17).Value = 570.8179886640253 ; $ws.Cells.Item(13, 18).Value = 5137.361897976228 ; $ws.Cells.Item(13, 19).Value = 0.02690672030782007 ; $ws.Cells.Item(13, 20).Value = 0.02797665847190844
$ws.Cells.Item(14, 7).Value = 3.072172 ; $ws.Cells.Item(14, 8).Value = 9.216516 ; $ws.Cells.Item(14, 9).Value = 0.07625183070469947 ; $ws.Cells.Item(14, 10).Value = 0.07781775827797052 ; $ws.Cells.Item(14, 13).Value = 137.0717086666666 ; $ws.Cells.Item(14, 14).Value = 411.2151259999999 ; $ws.Cells.Item(14, 15).Value = 0.2603191943704447 ; $ws.Cells.Item(14, 16).Value = 0.2652240042658267 ; $ws.Cells.Item(14, 17).Value = 421.1078653578906 ; $ws.Cells.Item(14, 18).Value = 3789.970788221015 ; $ws.Cells.Item(14, 19).Value = 0.0198498151383189 ; $ws.Cells.Item(14, 20).Value = 0.02063913745347352
$ws.Cells.Item(15, 7).Value = 3.072172 ; $ws.Cells.Item(15, 8).Value = 9.216516 ; $ws.Cells.Item(15, 9).Value = 0.07625183070469947 ; $ws.Cells.Item(15, 10).Value = 0.07781775827797052 ; $ws.Cells.Item(15, 13).Value = 29.2127365 ; $ws.Cells.Item(15, 14).Value = 58.425473 ; $ws.Cells.Item(15, 15).Value = 0.05547925319534149 ; $ws.Cells.Item(15, 16).Value = 0.03768304451958546 ; $ws.Cells.Item(15, 17).Value = 89.746551118678 ; $ws.Cells.Item(15, 18).Value = 538.479306712068 ; $ws.Cells.Item(15, 19).Value = 0.004230394622274336 ; $ws.Cells.Item(15, 20).Value = 0.002932410049603103
$ws.Cells.Item(16, 7).Value = 3.072172 ; $ws.Cells.Item(16, 8).Value = 9.216516 ; $ws.Cells.Item(16, 9).Value = 0.07625183070469947 ; $ws.Cells.Item(16, 10).Value = 0.07781775827797052 ; $ws.Cells.Item(16, 13).Value = 171.5584106666666 ; $ws.Cells.Item(16, 14).Value = 514.6752319999999 ; $ws.Cells.Item(16, 15).Value = 0.3258144783240821 ; $ws.Cells.Item(16, 16).Value = 0.331953319069988 ; $ws.Cells.Item(16, 17).Value = 527.0569456146346 ; $ws.Cells.Item(16, 18).Value = 4743.512510531712 ; $ws.Cells.Item(16, 19).Value = 0.02484395044230788 ; $ws.Cells.Item(16, 20).Value = 0.02583186314295835
$ws.Cells.Item(17, 7).Value = 2.432257 ; $ws.Cells.Item(17, 8).Value = 4.864514 ; $ws.Cells.Item(17, 9).Value = 0.06036903174507163 ; $ws.Cells.Item(17, 10).Value = 0.041072523998418 ; $ws.Cells.Item(17, 13).Value = 2.906846333333333 ; $ws.Cells.Item(17, 14).Value = 8.720538999999999 ; $ws.Cells.Item(17, 15).Value = 0.005520525738044089 ; $ws.Cells.Item(17, 16).Value = 0.005624540846623205 ; $ws.Cells.Item(17, 17).Value = 7.070197342174332 ; $ws.Cells.Item(17, 18).Value = 42.42118405304599 ; $ws.Cells.Item(17, 19).Value = 0.0003332687935294686 ; $ws.Cells.Item(17, 20).Value = 0.0002310140889030139
$ws.Cells.Item(18, 7).Value = 2.432257 ; $ws.Cells.Item(18, 8).Value = 4.864514 ; $ws.Cells.Item(18, 9).Value = 0.06036903174507163 ; $ws.Cells.Item(18, 10).Value = 0.041072523998418 ; $ws.Cells.Item(18, 13).Value = 185.8027443333333 ; $ws.Cells.Item(18, 14).Value = 557.408233 ; $ws.Cells.Item(18, 15).Value = 0.3528665483720876 ; $ws.Cells.Item(18, 16).Value = 0.3595150912979765 ; $ws.Cells.Item(18, 17).Value = 451.9200255239603 ; $ws.Cells.Item(18, 18).Value = 2711.520153143762 ; $ws.Cells.Item(18, 19).Value = 0.02130221186044842 ; $ws.Cells.Item(18, 20).Value = 0.01476619221512958
$ws.Cells.Item(19, 7).Value = 2.432257 ; $ws.Cells.Item(19, 8).Value = 4.864514 ; $ws.Cells.Item(19, 9).Value = 0.06036903174507163 ; $ws.Cells.Item(19, 10).Value = 0.041072523998418 ; $ws.Cells.Item(19, 13).Value = 137.0717086666666 ; $ws.Cells.Item(19, 14).Value = 411.2151259999999 ; $ws.Cells.Item(19, 15).Value = 0.2603191943704447 ; $ws.Cells.Item(19, 16).Value = 0.2652240042658267 ; $ws.Cells.Item(19, 17).Value = 333.3936229064606 ; $ws.Cells.Item(19, 18).Value = 2000.361737438764 ; $ws.Cells.Item(19, 19).Value = 0.01571521770880084 ; $ws.Cells.Item(19, 20).Value = 0.01089341928016469
$ws.Cells.Item(20, 7).Value = 2.432257 ; $ws.Cells.Item(20, 8).Value = 4.864514 ; $ws.Cells.Item(20, 9).Value = 0.06036903174507163 ; $ws.Cells.Item(20, 10).Value = 0.041072523998418 ; $ws.Cells.Item(20, 13).Value = 29.2127365 ; $ws.Cells.Item(20, 14).Value = 58.425473 ; $ws.Cells.Item(20, 15).Value = 0.05547925319534149 ; $ws.Cells.Item(20, 16).Value = 0.03768304451958546 ; $ws.Cells.Item(20, 17).Value = 71.05288284128049 ; $ws.Cells.Item(20, 18).Value = 284.211531365122 ; $ws.Cells.Item(20, 19).Value = 0.003349228797342437 ; $ws.Cells.Item(20, 20).Value = 0.001547737750364128
$ws.Cells.Item(21, 7).Value = 2.432257 ; $ws.Cells.Item(21, 8).Value = 4.864514 ; $ws.Cells.Item(21, 9).Value = 0.06036903174507163 ; $ws.Cells.Item(21, 10).Value = 0.041072523998418 ; $ws.Cells.Item(21, 13).Value = 171.5584106666666 ; $ws.Cells.Item(21, 14).Value = 514.6752319999999 ; $ws.Cells.Item(21, 15).Value = 0.3258144783240821 ; $ws.Cells.Item(21, 16).Value = 0.331953319069988 ; $ws.Cells.Item(21, 17).Value = 417.2741452528746 ; $ws.Cells.Item(21, 18).Value = 2503.644871517248 ; $ws.Cells.Item(21, 19).Value = 0.01966910458495047 ; $ws.Cells.Item(21, 20).Value = 0.01363416066385659
$ws.Cells.Item(22, 7).Value = 6.158463999999999 ; $ws.Cells.Item(22, 8).Value = 18.475392 ; $ws.Cells.Item(22, 9).Value = 0.1528541222070204 ; $ws.Cells.Item(22, 10).Value = 0.1559931745083229 ; $ws.Cells.Item(22, 13).Value = 2.906846333333333 ; $ws.Cells.Item(22, 14).Value = 8.720538999999999 ; $ws.Cells.Item(22, 15).Value = 0.005520525738044089 ; $ws.Cells.Item(22, 16).Value = 0.005624540846623205 ; $ws.Cells.Item(22, 17).Value = 17.90170849736533 ; $ws.Cells.Item(22, 18).Value = 161.115376476288 ; $ws.Cells.Item(22, 19).Value = 0.0008438351158099927 ; $ws.Cells.Item(22, 20).Value = 0.0008773899818164839
$ws.Cells.Item(23, 7).Value = 6.158463999999999 ; $ws.Cells.Item(23, 8).Value = 18.475392 ; $ws.Cells.Item(23, 9).Value = 0.1528541222070204 ; $ws.Cells.Item(23, 10).Value = 0.1559931745083229 ; $ws.Cells.Item(23, 13).Value = 185.8027443333333 ; $ws.Cells.Item(23, 14).Value = 557.408233 ; $ws.Cells.Item(23, 15).Value = 0.3528665483720876 ; $ws.Cells.Item(23, 16).Value = 0.3595150912979765 ; $ws.Cells.Item(23, 17).Value = 1144.259512078037 ; $ws.Cells.Item(23, 18).Value = 10298.33560870233 ; $ws.Cells.Item(23, 19).Value = 0.05393710650763656 ; $ws.Cells.Item(23, 20).Value = 0.05608190037522089
$ws.Cells.Item(24, 7).Value = 6.158463999999999 ; $ws.Cells.Item(24, 8).Value = 18.475392 ; $ws.Cells.Item(24, 9).Value = 0.1528541222070204 ; $ws.Cells.Item(24, 10).Value = 0.1559931745083229 ; $ws.Cells.Item(24, 13).Value = 137.0717086666666 ; $ws.Cells.Item(24, 14).Value = 411.2151259999999 ; $ws.Cells.Item(24, 15).Value = 0.2603191943704447 ; $ws.Cells.Item(24, 16).Value = 0.2652240042658267 ; $ws.Cells.Item(24, 17).Value = 844.1511832421544 ; $ws.Cells.Item(24, 18).Value = 7597.360649179391 ; $ws.Cells.Item(24, 19).Value = 0.03979086194913305 ; $ws.Cells.Item(24, 20).Value = 0.04137313438123528
$ws.Cells.Item(25, 7).Value = 6.158463999999999 ; $ws.Cells.Item(25, 8).Value = 18.475392 ; $ws.Cells.Item(25, 9).Value = 0.1528541222070204 ; $ws.Cells.Item(25, 10).Value = 0.1559931745083229 ; $ws.Cells.Item(25, 13).Value = 29.2127365 ; $ws.Cells.Item(25, 14).Value = 58.425473 ; $ws.Cells.Item(25, 15).Value = 0.05547925319534149 ; $ws.Cells.Item(25, 16).Value = 0.03768304451958546 ; $ws.Cells.Item(25, 17).Value = 179.905586076736 ; $ws.Cells.Item(25, 18).Value = 1079.433516460416 ; $ws.Cells.Item(25, 19).Value = 0.008480232547874954 ; $ws.Cells.Item(25, 20).Value = 0.005878297739748595
$ws.Cells.Item(26, 7).Value = 6.158463999999999 ; $ws.Cells.Item(26, 8).Value = 18.475392 ; $ws.Cells.Item(26, 9).Value = 0.1528541222070204 ; $ws.Cells.Item(26, 10).Value = 0.1559931745083229 ; $ws.Cells.Item(26, 13).Value = 171.5584106666666 ; $ws.Cells.Item(26, 14).Value = 514.6752319999999 ; $ws.Cells.Item(26, 15).Value = 0.3258144783240821 ; $ws.Cells.Item(26, 16).Value = 0.331953319069988 ; $ws.Cells.Item(26, 17).Value = 1056.536295987882 ; $ws.Cells.Item(26, 18).Value = 9508.826663890943 ; $ws.Cells.Item(26, 19).Value = 0.04980208608656585 ; $ws.Cells.Item(26, 20).Value = 0.05178245203030164
